$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new row of data needs to be inserted before current row 675, shifting the
# existing rows 675-716 down to 676-717. The new row shares the same date
# ("2026/01/18") and weekday ("day") text as row 674, so copy row 674 down
# (which preserves the existing text/shared-string cell formatting exactly)
# and then fix up the two numeric columns for the newly inserted row.
$ws.Rows.Item(674).Copy()
$ws.Rows.Item(675).Insert(-4121)

$ws.Range("C675").Value = 16
$ws.Range("D675").Value = 201
